$d = $word.ActiveDocument

# The document currently ends with:
#   ... <bookmark paragraph _GoBack> <trailing empty paragraph>
# We need to:
#  1. Insert a block of new paragraphs (ideas, headings, semaphore notes) right
#     before the bookmark paragraph.
#  2. Add a final "Signal" run inside the bookmark paragraph (before the bookmark).
#  3. Remove the trailing empty paragraph.

$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n - 1)

# 1) Insert a fresh empty paragraph right before the bookmark paragraph, then
#    replace its contents with the full OOXML block of new paragraphs.
$bookmarkPara.Range.InsertParagraphBefore()
$n = $d.Paragraphs.Count
$insertedPara = $d.Paragraphs.Item($n - 2)

$newParagraphsXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Idea 1; Funciona por turnos</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Idea 2: Miro y si no hay nadie levanto la mano.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Idea 3: Levanto la mano y luego miro</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Idea 4: Levanto la mano si no es mi turno, espero a mi turno y cuando toque levanto la mano y entro.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Sincronización alto nivel</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Semáforos:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Contador + lista de procesos en espera.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Iniciarl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(semáforo, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>valir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Wait</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (semáforo):  c&gt;</w:t></w:r><w:r><w:t xml:space="preserve">0 =&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Decrementa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> el contador /c==</w:t></w:r><w:r><w:t xml:space="preserve">0 suspende el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wait</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $newParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertedPara.Range.InsertXML($packageXml)

# 2) Re-find the bookmark paragraph (still the last-but-one paragraph) and add
#    the "Signal" run right before its bookmark.
$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n - 1)
$bookmarkPara.Range.InsertBefore("Signal")

# 3) Delete the trailing empty paragraph (merge its mark away) so the bookmark
#    paragraph becomes the last paragraph of the body again.
$n = $d.Paragraphs.Count
$trailing = $d.Paragraphs.Item($n)
$mergeRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$mergeRange.Delete()
